# BOT; UPDATE DATA
# Adds the 2020-05-06 (serial 43957) daily PCR/infection figures to the
# "all", "kobe" and "other" sheets, bumps the kobe H83 correction
# (8 -> 9), and refreshes the shared "municipal outsiders" footnote so it
# covers the two newly added case numbers (268, 272) and the updated
# count (12 -> 14 cases).

$wb = $excel.ActiveWorkbook

$footnote = "※　24・34・53・58・59・60・158・161・163・192・237・248・268・272例目（計14件）は市外在住者です。"

# ---------------------------------------------------------------------
# Sheet "all": insert a new row 29 (pushes the footnote row down to 30)
# and populate it with the new day's figures.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")

$wsAll.Rows("29").Insert()

$wsAll.Range("A29").Value = 43957
$wsAll.Range("B29").Value = 272
$wsAll.Range("C29").Value = 268
$wsAll.Range("D29").Value = 102
$wsAll.Range("E29").Value = 92
$wsAll.Range("F29").Value = 10
$wsAll.Range("G29").Value = 7
$wsAll.Range("H29").Value = 159

$wsAll.Range("B30").Value = $footnote

$wsAll.Activate()
$wsAll.Range("B30").Select()
$excel.ActiveWindow.ScrollRow = 28

# ---------------------------------------------------------------------
# Sheet "kobe": correct H83, fill in the already-present (but blank)
# row 84, and refresh the same shared footnote text on row 85.
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

$wsKobe.Range("H83").Value = 9

$wsKobe.Range("A84").Value = 43957
$wsKobe.Range("B84").Value = 72
$wsKobe.Range("C84").Value = 2351
$wsKobe.Range("D84").Value = 2
$wsKobe.Range("E84").Value = 272
$wsKobe.Range("F84").Value = 97
$wsKobe.Range("G84").Value = 88
$wsKobe.Range("H84").Value = 9
$wsKobe.Range("I84").Value = 7
$wsKobe.Range("J84").Value = 152

$wsKobe.Range("B85").Value = $footnote

$wsKobe.Activate()
$wsKobe.Range("B86").Select()
$excel.ActiveWindow.ScrollRow = 82

# ---------------------------------------------------------------------
# Sheet "other": fill in the already-present (but blank) row 59; copy
# the formatting down from row 58 first so the per-column styles match,
# then overwrite with the real values.
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")

$wsOther.Range("A58:H58").Copy()
$wsOther.Range("A59:H59").PasteSpecial(-4122)

$wsOther.Range("A59").Value = 43957
$wsOther.Range("B59").Value = 0
$wsOther.Range("C59").Value = 12
$wsOther.Range("D59").Value = 5
$wsOther.Range("E59").Value = 4
$wsOther.Range("F59").Value = 1
$wsOther.Range("G59").Value = 0
$wsOther.Range("H59").Value = 7

$wsOther.Activate()
$wsOther.Range("J59").Select()
$excel.ActiveWindow.ScrollRow = 56
$excel.ActiveWindow.ScrollColumn = 5
